$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$notesPage = $s.NotesPage
$shp = $notesPage.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The notes placeholder holds 4 paragraphs separated by "`r" when read back.
# Paragraph 3 currently ends in two runs ("... felhasználó " + "alatt is.");
# merge them into a single run reading "... felhasználó alatt is" (no
# trailing period).
$fullText = $tr.Text
$paragraphs = $fullText -split "`r"

$target = "Növekvő adatok kezelése könnyebb, kiváló teljesítmény sok felhasználó alatt is."
$replacement = "Növekvő adatok kezelése könnyebb, kiváló teljesítmény sok felhasználó alatt is"

for ($i = 0; $i -lt $paragraphs.Count; $i++) {
    if ($paragraphs[$i] -eq $target) {
        $paragraphs[$i] = $replacement
    }
}

$tr.Text = [string]::Join([string][char]10, $paragraphs)
